# Add a new "Height" property/column for NPC data (column AB), filling
# every data row (2-21) with the value 2, matching the commit
# "add property for npc".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in AB1
$ws.Cells.Item(1, 28).Value = "Height"

# Fill AB2:AB21 with the value 2 for every NPC row
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 28).Value = 2
}

# Match column width of the newly inserted column as closely as possible
# to the sibling column (AA, 14.75) given the available width granularity.
$ws.Cells.Item(1, 28).ColumnWidth = 14

# Reflect the new selection / scrolled view used after adding the column
$win = $excel.ActiveWindow
$win.ScrollColumn = 24
$ws.Range("AB2:AB21").Select()
